# shop_item.xlsx edit:
# "Add showing how full shelf is, Next add the same for the other stands"
#
# For the first "stand" of items (rows 2-5, already unlocked) we add back a
# sprite_path value (now just the short sprite key instead of a full
# res:// path) and drop the now-unused average_color column (L).
# For the next batch of items (rows 6-12) we unlock them (G=1) and give them
# a sprite_path value too, same as the first stand.
# Rows 13-21 (further stands) are left untouched for now, matching the
# "Next add the same for the other stands" follow-up note.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sprite_path (column H) for the already-unlocked first stand ---------
# These replace the old full "res://Asset/ShopItem/XXX.png" style strings
# that used to live in column H with short sprite keys.
$ws.Range("H2").Value = "Bread"
$ws.Range("H3").Value = "Mead"
$ws.Range("H4").Value = "Apple"
$ws.Range("H5").Value = "Cheese"

# --- unlock the next stand of items and give them sprite_path values -----
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = "BootsLeather"

$ws.Range("G7").Value = 1
$ws.Range("H7").Value = "SwordIron"

$ws.Range("G8").Value = 1
$ws.Range("H8").Value = "PotionWideGreen"

$ws.Range("G9").Value = 1
$ws.Range("H9").Value = "WoolWhite"

$ws.Range("G10").Value = 1
$ws.Range("H10").Value = "ShieldWood"

$ws.Range("G11").Value = 1
$ws.Range("H11").Value = "Honeycomb"

$ws.Range("G12").Value = 1
$ws.Range("H12").Value = "MugClay"

# --- drop the now-unused average_color column (L) ------------------------
$ws.Range("L1:L5").ClearContents()

# --- cosmetic sheet-view tweaks ------------------------------------------
$ws.Range("H1").ColumnWidth = 18.15
$ws.Range("G24").Select()
